$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets("ALC")
$ws.Range("H106").Value = 45979400
$ws.Range("I106").Value = 25643582
$ws.Range("J106").Value = 62502250
$ws.Range("K106").Value = 25643582
$ws.Range("L106").Value = 62502250
$ws.Range("M106").Value = -25642951
$ws.Range("N106").Value = -62503512

$ws.Range("H116").Value = 7129.5713
$ws.Range("I116").Value = 9401.154
$ws.Range("J116").Value = 3438.25
$ws.Range("K116").Value = 9401.154
$ws.Range("L116").Value = 3438.25
$ws.Range("M116").Value = -5959.154
$ws.Range("N116").Value = -10322.25

$ws.Range("H132").Value = 6173668
$ws.Range("I132").Value = 628.0465
$ws.Range("J132").Value = 30304642
$ws.Range("K132").Value = 1884.1395
$ws.Range("L132").Value = 90913926
$ws.Range("M132").Value = 645.8604999999998
$ws.Range("N132").Value = -90918986

$ws.Range("H135").Value = 1374.4
$ws.Range("I135").Value = 1101.8667
$ws.Range("J135").Value = 3827.2
$ws.Range("K135").Value = 9916.800300000001
$ws.Range("L135").Value = 34444.8
$ws.Range("M135").Value = -7381.800300000001
$ws.Range("N135").Value = -39514.8

$ws.Range("H137").Value = 1285.7273
$ws.Range("I137").Value = 1039.2778
$ws.Range("J137").Value = 1752.6842
$ws.Range("K137").Value = 3117.8334
$ws.Range("L137").Value = 5258.0526
$ws.Range("M137").Value = -567.8334000000004
$ws.Range("N137").Value = -10358.0526

$ws.Range("H138").Value = 2349.7896
$ws.Range("I138").Value = 827.7778
$ws.Range("J138").Value = 4559.161
$ws.Range("K138").Value = 2483.3334
$ws.Range("L138").Value = 13677.483
$ws.Range("M138").Value = 2656.6666
$ws.Range("N138").Value = -23957.483

$ws = $wb.Sheets("ARM")
$ws.Range("H32").Value = 5641.925
$ws.Range("I32").Value = 4242.6094
$ws.Range("J32").Value = 11239.1875
$ws.Range("K32").Value = 4242.6094
$ws.Range("L32").Value = 11239.1875
$ws.Range("M32").Value = -3955.6094
$ws.Range("N32").Value = -11813.1875

$ws.Range("H45").Value = 7268.353
$ws.Range("I45").Value = 10697.2
$ws.Range("J45").Value = 2370
$ws.Range("K45").Value = 10697.2
$ws.Range("L45").Value = 2370
$ws.Range("M45").Value = -10320.2
$ws.Range("N45").Value = -3124

$ws.Range("H61").Value = 3677.1904
$ws.Range("I61").Value = 3711.05
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 3711.05
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -3499.05
$ws.Range("N61").Value = -3424

$ws.Range("H74").Value = 13514683
$ws.Range("I74").Value = 928.4666999999999
$ws.Range("J74").Value = 71430776
$ws.Range("K74").Value = 928.4666999999999
$ws.Range("L74").Value = 71430776
$ws.Range("M74").Value = -54.46669999999995
$ws.Range("N74").Value = -71432524

$ws.Range("H77").Value = 13514683
$ws.Range("I77").Value = 928.4666999999999
$ws.Range("J77").Value = 71430776
$ws.Range("K77").Value = 4642.3335
$ws.Range("L77").Value = 357153880
$ws.Range("M77").Value = -274.3334999999997
$ws.Range("N77").Value = -357162616

$ws.Range("H123").Value = 29714.5
$ws.Range("J123").Value = 29714.5
$ws.Range("L123").Value = 29714.5
$ws.Range("N123").Value = -39514.5

$ws.Range("H125").Value = 52293.5
$ws.Range("J125").Value = 52293.5
$ws.Range("L125").Value = 52293.5
$ws.Range("N125").Value = -62133.5

$ws.Range("H132").Value = 2911.0232
$ws.Range("I132").Value = 2141.8064
$ws.Range("J132").Value = 4898.1665
$ws.Range("K132").Value = 6425.4192
$ws.Range("L132").Value = 14694.4995
$ws.Range("M132").Value = -3895.4192
$ws.Range("N132").Value = -19754.4995

$ws.Range("H136").Value = 3677.1904
$ws.Range("I136").Value = 3711.05
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 11133.15
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -8583.150000000001
$ws.Range("N136").Value = -14100

$ws = $wb.Sheets("BSM")
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("K16").Value = 5000
$ws.Range("M16").Value = -4830

$ws.Range("H86").Value = 1785.7142
$ws.Range("I86").Value = 1744.4445
$ws.Range("J86").Value = 1860
$ws.Range("K86").Value = 1744.4445
$ws.Range("L86").Value = 1860
$ws.Range("M86").Value = -621.4445000000001
$ws.Range("N86").Value = -4106

$ws.Range("H89").Value = 1785.7142
$ws.Range("I89").Value = 1744.4445
$ws.Range("J89").Value = 1860
$ws.Range("K89").Value = 8722.2225
$ws.Range("L89").Value = 9300
$ws.Range("M89").Value = -3106.2225
$ws.Range("N89").Value = -20532

$ws.Range("H99").Value = 166668530
$ws.Range("I99").Value = 200001250
$ws.Range("K99").Value = 200001250
$ws.Range("M99").Value = -199999752

$ws.Range("H107").Value = 166668580
$ws.Range("I107").Value = 200002000
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 200002000
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -200000080
$ws.Range("N107").Value = -5340

$ws.Range("H134").Value = 4211.0425
$ws.Range("I134").Value = 5029.6562
$ws.Range("K134").Value = 15088.9686
$ws.Range("M134").Value = -12553.9686

$ws = $wb.Sheets("CRP")
$ws.Range("H12").Value = 1805.3
$ws.Range("I12").Value = 1336.1428
$ws.Range("J12").Value = 2900
$ws.Range("K12").Value = 1336.1428
$ws.Range("L12").Value = 2900
$ws.Range("M12").Value = -1166.1428
$ws.Range("N12").Value = -3240

$ws.Range("H31").Value = 6611404
$ws.Range("I31").Value = 1717.8334
$ws.Range("J31").Value = 14322704
$ws.Range("K31").Value = 1717.8334
$ws.Range("L31").Value = 14322704
$ws.Range("M31").Value = -1422.8334
$ws.Range("N31").Value = -14323294

$ws.Range("H34").Value = 6611404
$ws.Range("I34").Value = 1717.8334
$ws.Range("J34").Value = 14322704
$ws.Range("K34").Value = 1717.8334
$ws.Range("L34").Value = 14322704
$ws.Range("M34").Value = -1515.8334
$ws.Range("N34").Value = -14323108

$ws.Range("H99").Value = 10428609
$ws.Range("I99").Value = 13630.286
$ws.Range("K99").Value = 13630.286
$ws.Range("M99").Value = -12132.286

$ws.Range("H126").Value = 10428609
$ws.Range("I126").Value = 13630.286
$ws.Range("K126").Value = 40890.858
$ws.Range("M126").Value = -38420.858

$ws.Range("H132").Value = 2326788
$ws.Range("I132").Value = 2899566
$ws.Range("J132").Value = 1983.3529
$ws.Range("K132").Value = 8698698
$ws.Range("L132").Value = 5950.0587
$ws.Range("M132").Value = -8696168
$ws.Range("N132").Value = -11010.0587

$ws.Range("H134").Value = 4168029.5
$ws.Range("I134").Value = 6668262
$ws.Range("J134").Value = 975.3333
$ws.Range("K134").Value = 20004786
$ws.Range("L134").Value = 2925.9999
$ws.Range("M134").Value = -20002251
$ws.Range("N134").Value = -7995.9999

$ws = $wb.Sheets("CUL")
$ws.Range("H113").Value = 2500564.5
$ws.Range("I113").Value = 3846750.8
$ws.Range("K113").Value = 11540252.4
$ws.Range("M113").Value = -11538082.4

$ws = $wb.Sheets("GSM")
$ws.Range("H12").Value = 13437.5
$ws.Range("J12").Value = 12983.333
$ws.Range("L12").Value = 12983.333
$ws.Range("N12").Value = -13263.333

$ws.Range("H80").Value = 3433.75
$ws.Range("I80").Value = 2490
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2490
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1492
$ws.Range("N80").Value = -5996

$ws.Range("H83").Value = 3433.75
$ws.Range("I83").Value = 2490
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 12450
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -7458
$ws.Range("N83").Value = -29984

$ws.Range("H97").Value = 738.53845
$ws.Range("I97").Value = 623.75
$ws.Range("J97").Value = 922.2
$ws.Range("K97").Value = 623.75
$ws.Range("L97").Value = 922.2
$ws.Range("M97").Value = -127.75
$ws.Range("N97").Value = -1914.2

$ws.Range("H123").Value = 23230.62
$ws.Range("J123").Value = 23230.62
$ws.Range("L123").Value = 23230.62
$ws.Range("N123").Value = -28130.62

$ws.Range("H132").Value = 3789662
$ws.Range("I132").Value = 4506194.5
$ws.Range("J132").Value = 2277.2856
$ws.Range("K132").Value = 13518583.5
$ws.Range("L132").Value = 6831.8568
$ws.Range("M132").Value = -13516053.5
$ws.Range("N132").Value = -11891.8568

$ws = $wb.Sheets("LTW")
$ws.Range("H46").Value = 13889464
$ws.Range("I46").Value = 41667132
$ws.Range("J46").Value = 629.375
$ws.Range("K46").Value = 41667132
$ws.Range("L46").Value = 629.375
$ws.Range("M46").Value = -41666944
$ws.Range("N46").Value = -1005.375

$ws.Range("H122").Value = 6793697.5
$ws.Range("I122").Value = 10215653
$ws.Range("J122").Value = 2002960
$ws.Range("K122").Value = 30646959
$ws.Range("L122").Value = 6008880
$ws.Range("M122").Value = -30644509
$ws.Range("N122").Value = -6013780

$ws.Range("H132").Value = 16671869
$ws.Range("I132").Value = 19053030
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 57159090
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = -57156560
$ws.Range("N132").Value = -16310

$ws.Range("H136").Value = 9584.091
$ws.Range("I136").Value = 7174.25
$ws.Range("J136").Value = 16010.333
$ws.Range("K136").Value = 21522.75
$ws.Range("L136").Value = 48030.999
$ws.Range("M136").Value = -18972.75
$ws.Range("N136").Value = -53130.999

$ws = $wb.Sheets("WVR")
$ws.Range("H14").Value = 9288.360000000001
$ws.Range("I14").Value = 9175
$ws.Range("J14").Value = 9341.706
$ws.Range("K14").Value = 9175
$ws.Range("L14").Value = 9341.706
$ws.Range("M14").Value = -9007
$ws.Range("N14").Value = -9677.706

$ws.Range("H122").Value = 1114.75
$ws.Range("I122").Value = 1157.2858
$ws.Range("J122").Value = 987.1429000000001
$ws.Range("K122").Value = 3471.8574
$ws.Range("L122").Value = 2961.4287
$ws.Range("M122").Value = -1021.8574
$ws.Range("N122").Value = -7861.4287

$ws.Range("H123").Value = 40428
$ws.Range("J123").Value = 40428
$ws.Range("L123").Value = 40428
$ws.Range("N123").Value = -50228

$ws.Range("H126").Value = 1377.6364
$ws.Range("I126").Value = 864.6667
$ws.Range("J126").Value = 1570
$ws.Range("K126").Value = 2594.0001
$ws.Range("L126").Value = 4710
$ws.Range("M126").Value = -124.0001000000002
$ws.Range("N126").Value = -9650

$ws.Range("H132").Value = 1149.5714
$ws.Range("I132").Value = 880.8461
$ws.Range("J132").Value = 1586.25
$ws.Range("K132").Value = 2642.5383
$ws.Range("L132").Value = 4758.75
$ws.Range("M132").Value = -112.5383000000002
$ws.Range("N132").Value = -9818.75

$ws.Range("H136").Value = 13074745
$ws.Range("I136").Value = 3440.9
$ws.Range("J136").Value = 31748038
$ws.Range("K136").Value = 10322.7
$ws.Range("L136").Value = 95244114
$ws.Range("M136").Value = -7772.700000000001
$ws.Range("N136").Value = -95249214
